$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.440.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.541.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.93"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.28"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.943.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.544.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.365.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0941"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "284.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.07%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.79"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.54"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.36"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0782"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.67"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.34"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0310"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.009.77"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.12"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.29"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.64%  "
